$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.623.53"
$ws.Range("E2").Value = "  +0.44%  "

$ws.Range("D3").Value = "1.925.54"
$ws.Range("E3").Value = "  +0.13%  "

$ws.Range("E4").Value = "  +0.79%  "

$ws.Range("D5").Value = "'327.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.42%  "

$ws.Range("E6").Value = "  +0.73%  "

$ws.Range("D7").Value = "'0.4819"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.51%  "

$ws.Range("D8").Value = "'0.4055"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.85%  "

$ws.Range("D9").Value = "'0.08192"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.12%  "

$ws.Range("E10").Value = "  -1.88%  "

$ws.Range("E11").Value = "  -0.44%  "

$ws.Range("D12").Value = "1.901.57"
$ws.Range("E12").Value = "  +0.00%  "

$ws.Range("D13").Value = "'6.064"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.16%  "

$ws.Range("D14").Value = "'7.282"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.65%  "

$ws.Range("D15").Value = "'91.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.19%  "

$ws.Range("D16").Value = "'0.06869"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.53%  "

$ws.Range("D17").Value = "'1.014"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.80%  "

$ws.Range("E18").Value = "  -0.23%  "

$ws.Range("D19").Value = "'17.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.08%  "

$ws.Range("E20").Value = "  +0.63%  "

$ws.Range("D21").Value = "29.606.28"
$ws.Range("E21").Value = "  +0.29%  "

$ws.Range("D22").Value = "'5.646"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").Value = "'11.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.58%  "

$ws.Range("D24").Value = "'2.202"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.03%  "

$ws.Range("D25").Value = "2.090.08"
$ws.Range("E25").Value = "  -2.10%  "

$ws.Range("E26").Value = "  -0.16%  "

$ws.Range("D27").Value = "'6.374"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.05%  "

$ws.Range("D28").Value = "'19.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.53%  "

$ws.Range("D29").Value = "'2.084"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.08%  "

$ws.Range("D30").Value = "'120.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.06%  "

$ws.Range("D31").Value = "'1.005"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.23%  "

$ws.Range("D32").Value = "'0.09597"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.23%  "

$ws.Range("D33").Value = "'5.609"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.19%  "

$ws.Range("D34").Value = "'3.559"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.34%  "

$ws.Range("D35").Value = "'1.391"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.39%  "

$ws.Range("D36").Value = "'0.06498"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.59%  "

$ws.Range("D38").Value = "'1.208"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.08%  "

$ws.Range("D39").Value = "'0.5922"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.20%  "

$ws.Range("D40").Value = "'10.75"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.92%  "

$ws.Range("D41").Value = "'7.853"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.04%  "

$ws.Range("D42").Value = "'2.526"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.29%  "

$ws.Range("D43").Value = "'0.1841"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.35%  "

$ws.Range("D44").Value = "'1.281"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.02%  "

$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "'0.07524"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.33%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'12.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.10%  "

$ws.Range("D47").Value = "'0.5542"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.02%  "

$ws.Range("D48").Value = "'1.958"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.40%  "

$ws.Range("D49").Value = "'117.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.31%  "

$ws.Range("D50").Value = "'2.425"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.47%  "

$ws.Range("D51").Value = "'71.85"
$ws.Range("D51").Style = "Normal"
